$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: clone current rows 234 and 235 (A:R) into new rows 236 and 237
$ws.Range("A234:R234").Copy($ws.Range("A236:R236"))
$ws.Range("A235:R235").Copy($ws.Range("A237:R237"))

# Step 2: capture original D/J/K/L/M/O/P for rows 168..235 before overwriting
$origD = @{}
$origJ = @{}
$origK = @{}
$origL = @{}
$origM = @{}
$origO = @{}
$origP = @{}
for ($r = 168; $r -le 235; $r++) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value2
    $origJ[$r] = $ws.Cells.Item($r, 10).Value2
    $origK[$r] = $ws.Cells.Item($r, 11).Value2
    $origL[$r] = $ws.Cells.Item($r, 12).Value2
    $origM[$r] = $ws.Cells.Item($r, 13).Value2
    $origO[$r] = $ws.Cells.Item($r, 15).Value2
    $origP[$r] = $ws.Cells.Item($r, 16).Value2
}

# Step 3: shift rows 170..235 to take the (D,J,K,L,M,O,P) values from two rows above
for ($r = 235; $r -ge 170; $r--) {
    $src = $r - 2
    $ws.Cells.Item($r, 4).Value = $origD[$src]
    $ws.Cells.Item($r, 10).Value = $origJ[$src]
    $ws.Cells.Item($r, 11).Value = $origK[$src]
    $ws.Cells.Item($r, 12).Value = $origL[$src]
    $ws.Cells.Item($r, 13).Value = $origM[$src]
    $ws.Cells.Item($r, 15).Value = $origO[$src]
    $ws.Cells.Item($r, 16).Value = $origP[$src]
}

# Step 4: rows 168 and 169 get the brand-new week of data (J unchanged, D/K/L/M/P updated)
$ws.Cells.Item(168, 4).Value = 45141
$ws.Cells.Item(168, 11).Value = 600
$ws.Cells.Item(168, 12).Value = 700
$ws.Cells.Item(168, 13).Value = 650
$ws.Cells.Item(168, 16).Value = 650

$ws.Cells.Item(169, 4).Value = 45141
$ws.Cells.Item(169, 11).Value = 500
$ws.Cells.Item(169, 12).Value = 500
$ws.Cells.Item(169, 13).Value = 500
$ws.Cells.Item(169, 16).Value = 500
